$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting for the Price (D) and Volume(1h) (E) columns so that
# numeric-looking strings (e.g. "212.26") are stored as text, matching the
# original inline-string cell contents instead of being parsed as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.945.09"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "1.637.00"
$ws.Range("E3").Value = "  -0.62%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "212.26"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "23.30"
$ws.Range("E8").Value = "  -1.41%  "
$ws.Range("E9").Value = "  -2.33%  "
$ws.Range("E10").Value = "  -0.18%  "
$ws.Range("E11").Value = "  +1.13%  "
$ws.Range("D12").Value = "1.870.19"
$ws.Range("E12").Value = "  -0.60%  "
$ws.Range("D13").Value = "1.640.09"
$ws.Range("E13").Value = "  -0.44%  "
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("D15").Value = "0.568"
$ws.Range("E15").Value = "  +0.72%  "
$ws.Range("D16").Value = "65.30"
$ws.Range("E16").Value = "  -0.66%  "
$ws.Range("D17").Value = "27.955.71"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").Value = "231.04"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").Value = "0.0₃0720"
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("D20").Value = "7.53"
$ws.Range("E20").Value = "  -1.82%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").Value = "10.40"
$ws.Range("E22").Value = "  -2.68%  "
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("D24").Value = "2.07"
$ws.Range("E24").Value = "  -3.89%  "
$ws.Range("D25").Value = "153.56"
$ws.Range("E25").Value = "  +1.23%  "
$ws.Range("D26").Value = "6.98"
$ws.Range("E26").Value = "  +0.92%  "
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("E28").Value = "  -0.71%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  -0.46%  "
$ws.Range("D31").Value = "0.0482"
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("E32").Value = "  +1.26%  "
$ws.Range("D33").Value = "1.410.76"
$ws.Range("E33").Value = "  -3.02%  "
$ws.Range("D34").Value = "3.07"
$ws.Range("E34").Value = "  -1.43%  "
$ws.Range("E35").Value = "  +1.29%  "
$ws.Range("E36").Value = "  +1.54%  "
$ws.Range("E37").Value = "  +0.45%  "
$ws.Range("D38").Value = "0.562"
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("D39").Value = "0.929"
$ws.Range("E39").Value = "  +1.39%  "
$ws.Range("D40").Value = "0.876"
$ws.Range("E40").Value = "  -1.34%  "
$ws.Range("E41").Value = "  +0.66%  "
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").Value = "67.04"
$ws.Range("E43").Value = "  -3.45%  "
$ws.Range("D44").Value = "5.52"
$ws.Range("E44").Value = "  +2.61%  "
$ws.Range("E45").Value = "  +1.77%  "
$ws.Range("D46").Value = "2.20"
$ws.Range("E46").Value = "  -1.22%  "
$ws.Range("D47").Value = "1.779.30"
$ws.Range("E47").Value = "  -0.67%  "
$ws.Range("D48").Value = "87.93"
$ws.Range("E48").Value = "  -1.15%  "
$ws.Range("E49").Value = "  -0.45%  "
$ws.Range("E50").Value = "  -0.46%  "
$ws.Range("E51").Value = "  -2.23%  "
